$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix duplicate company names in column B (checked & renamed the repeats) ---
$ws.Range("B15").Value = "NCR-2"
$ws.Range("B28").Value = "TCP Wave-2"
$ws.Range("B32").Value = "Zensar Technologies-2"
$ws.Range("B44").Value = "NCR-3"
$ws.Range("B42").Value = "Netcracker-3"
$ws.Range("B39").Value = "Virtusa-2"
$ws.Range("B75").Value = "Hexaware-2"

# --- Highlight duplicate company names in column B ---
$dupRng = $ws.Range("B1:B1048576")
$fc = $dupRng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# --- Turn on AutoFilter for the SALARY column ---
$filterRng = $ws.Range("C1:C111")
$filterRng.AutoFilter() | Out-Null

# --- Sort the data range by SALARY (column C) descending ---
$fullRng = $ws.Range("A1:D111")
$sortKey = $ws.Range("C1:C111")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($fullRng)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Record the (hidden) filter database name Excel normally stamps on AutoFilter ---
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$C`$1:`$C`$111")
$fdName.Visible = $false

# --- Leave the selection where the user left it after the sort/filter pass ---
$ws.Range("A2:A111").Select() | Out-Null
